$d = $word.ActiveDocument

$replacements = @(
    @("84×12=", "22×29="),
    @("65×23=", "51×54="),
    @("12×59=", "22×71="),
    @("49×81=", "15×75="),
    @("48×58=", "49×75="),
    @("29×43=", "46×79="),
    @("34×46=", "34×91="),
    @("95×53=", "57×77="),
    @("53×74=", "58×59="),
    @("47×62=", "87×72="),
    @("33×45=", "36×76="),
    @("30×52=", "37×16="),
    @("91×84=", "80×79="),
    @("67×66=", "51×14="),
    @("18×89=", "82×93="),
    @("47×58=", "86×37="),
    @("87×71=", "20×37="),
    @("65×64=", "87×63="),
    @("22×96=", "43×74="),
    @("81×75=", "62×12="),
    @("25×53=", "36×53="),
    @("36×63=", "23×95="),
    @("25×48=", "87×83="),
    @("42×37=", "34×94="),
    @("16×32=", "88×18=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
